$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.675.70'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '1.808.45'
$ws.Range('D4').Value = '''0.999'
$ws.Range('D5').Value = '''225.55'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '''0.553'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D8').Value = '''32.71'
$ws.Range('E8').Value = '  +4.67%  '
$ws.Range('D9').Value = '''0.291'
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').Value = '''0.0713'
$ws.Range('E10').Value = '  +7.85%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '2.068.43'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '''11.13'
$ws.Range('E13').Value = '  -3.35%  '
$ws.Range('D14').Value = '1.807.11'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').Value = '34.688.22'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '''69.71'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '''254.57'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').Value = '0.0₃0805'
$ws.Range('E20').Value = '  +8.35%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '''10.94'
$ws.Range('E21').Value = '  +4.51%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '''4.26'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('D25').Value = '''161.74'
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('D26').Value = '''16.51'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '''0.0534'
$ws.Range('E30').Value = '  +3.38%  '
$ws.Range('D31').Value = '''3.81'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '''3.65'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('D35').Value = '1.440.29'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '''0.645'
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('E38').Value = '  +3.15%  '
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('D40').Value = '''0.959'
$ws.Range('E40').Value = '  +6.24%  '
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('E43').Value = '  +2.96%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''6.05'
$ws.Range('E44').Value = '  +5.60%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '''1.06'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = '''0.0494'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.963.39'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''106.33'
$ws.Range('E48').Value = '  +8.84%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''12.23'
$ws.Range('E49').Value = '  +3.06%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0127'
$ws.Range('E51').Value = '  +9.81%  '
